$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the data rows (2-8) with the new charge-type rows/labels ---

# Row 2: ItemStorage, 0-30 CMB/Day, 0 EUR
$ws.Cells.Item(2,1).Value = "ItemStorage"
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 30
$ws.Cells.Item(2,4).Value = "CMB/Day"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = "EUR"

# Row 3: ItemStorage, 31-60 CMB/Day, 0.35 EUR
$ws.Cells.Item(3,1).Value = "ItemStorage"
$ws.Cells.Item(3,2).Value = 31
$ws.Cells.Item(3,3).Value = 60
$ws.Cells.Item(3,4).Value = "CMB/Day"
$ws.Cells.Item(3,5).Value = 0.35
$ws.Cells.Item(3,6).Value = "EUR"

# Row 4: ItemStorage, 61-999 CMB/Day, 0.6 EUR
$ws.Cells.Item(4,1).Value = "ItemStorage"
$ws.Cells.Item(4,2).Value = 61
$ws.Cells.Item(4,3).Value = 999
$ws.Cells.Item(4,4).Value = "CMB/Day"
$ws.Cells.Item(4,5).Value = 0.6
$ws.Cells.Item(4,6).Value = "EUR"

# Row 5: OutboundHandling, 0-1 KG, 0.06 EUR
$ws.Cells.Item(5,1).Value = "OutboundHandling"
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 1
$ws.Cells.Item(5,4).Value = "KG"
$ws.Cells.Item(5,5).Value = 0.06
$ws.Cells.Item(5,6).Value = "EUR"

# Row 6: OutboundHandling, 1-10 KG, 0.17 EUR
$ws.Cells.Item(6,1).Value = "OutboundHandling"
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = 10
$ws.Cells.Item(6,4).Value = "KG"
$ws.Cells.Item(6,5).Value = 0.17
$ws.Cells.Item(6,6).Value = "EUR"

# Row 7: OutboundHandling, 1-10 KG, 0.17 EUR, over-limit 1 @ 0.15, 10-20-30
$ws.Cells.Item(7,1).Value = "OutboundHandling"
$ws.Cells.Item(7,2).Value = 1
$ws.Cells.Item(7,3).Value = 10
$ws.Cells.Item(7,4).Value = "KG"
$ws.Cells.Item(7,5).Value = 0.17
$ws.Cells.Item(7,6).Value = "EUR"
$ws.Cells.Item(7,7).Value = 1
$ws.Cells.Item(7,8).Value = 0.15
$ws.Cells.Item(7,9).Value = 10
$ws.Cells.Item(7,10).Value = 20
$ws.Cells.Item(7,11).Value = 30

# Row 8: OutboundHandling, 1-10 KG, 0.17 EUR, over-limit 1 @ 0.15, 20-30-40
$ws.Cells.Item(8,1).Value = "OutboundHandling"
$ws.Cells.Item(8,2).Value = 1
$ws.Cells.Item(8,3).Value = 10
$ws.Cells.Item(8,4).Value = "KG"
$ws.Cells.Item(8,5).Value = 0.17
$ws.Cells.Item(8,6).Value = "EUR"
$ws.Cells.Item(8,7).Value = 1
$ws.Cells.Item(8,8).Value = 0.15
$ws.Cells.Item(8,9).Value = 20
$ws.Cells.Item(8,10).Value = 30
$ws.Cells.Item(8,11).Value = 40

# --- Rows 9-11: drop the old data but keep the styled/empty currency cell ---
$ws.Range("A9:E9").Clear() | Out-Null
$ws.Range("F9").ClearContents() | Out-Null

$ws.Range("A10:E10").Clear() | Out-Null
$ws.Range("F10").ClearContents() | Out-Null

$ws.Range("A11:E11").Clear() | Out-Null
$ws.Range("F11").ClearContents() | Out-Null
$ws.Range("G11:J11").Clear() | Out-Null
$ws.Range("K11").ClearContents() | Out-Null

# --- Rows 12-15: drop all the old sample data entirely ---
$ws.Range("A12:J12").Clear() | Out-Null
$ws.Range("K12").ClearContents() | Out-Null

$ws.Range("A13:J13").Clear() | Out-Null
$ws.Range("K13").ClearContents() | Out-Null

$ws.Range("A14:J14").Clear() | Out-Null
$ws.Range("K14").ClearContents() | Out-Null

$ws.Range("A15:J15").Clear() | Out-Null
$ws.Range("K15").ClearContents() | Out-Null

# --- Drop the now-unused trailing placeholder rows 30-33 ---
$ws.Range("K30:K33").EntireRow.Delete() | Out-Null

# --- Restore the view: scroll back to A1 and select A3:A4 ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A3:A4").Select() | Out-Null
